# wms_kaart_database.xlsx — "Wijziging van kolomnamen uppercase en lowercase en testen"
#
# 1) Rename several column-name / code strings (uppercase/lowercase, dash -> underscore)
# 2) Re-apply the AutoFilter bookkeeping (adds the extra _FilterDatabase_0 / _0_0 names)
# 3) Update the saved view state (freeze-pane scroll position + selection) on Blad1
# 4) Small row-height tweak on header row 2
# 5) Column width adjustments on Blad1 (A:J)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Text replacements (these are shared strings, so every cell that used the
#    old text must be updated to the new text).
# ---------------------------------------------------------------------------

# "lokaalid, bgttype, plustype" -> "identificatie_lokaalid, bgt_type, plus_type"
$ws.Range("F2:F53").Value = "identificatie_lokaalid, bgt_type, plus_type"
$ws.Range("F69:F71").Value = "identificatie_lokaalid, bgt_type, plus_type"
$ws.Range("F73:F74").Value = "identificatie_lokaalid, bgt_type, plus_type"
$ws.Range("F76:F86").Value = "identificatie_lokaalid, bgt_type, plus_type"
$ws.Range("F95").Value = "identificatie_lokaalid, bgt_type, plus_type"
$ws.Range("F99:F108").Value = "identificatie_lokaalid, bgt_type, plus_type"
$ws.Range("F112:F115").Value = "identificatie_lokaalid, bgt_type, plus_type"

# "BGTPLUS_KST_cai-kast" -> "BGTPLUS_KST_cai_kast"
$ws.Range("B13").Value = "BGTPLUS_KST_cai_kast"

# "BGTPLUS_PUT_brandkraan_-put" -> "BGTPLUS_PUT_brandkraan__put"
$ws.Range("B29").Value = "BGTPLUS_PUT_brandkraan__put"

# "BGTPLUS_PUT_inspectie-_rioolput" -> "BGTPLUS_PUT_inspectie__rioolput"
$ws.Range("B30").Value = "BGTPLUS_PUT_inspectie__rioolput"

# "lokaalid, bgtfysvkn, plusfyskvkn" -> "identificatie_lokaalid, bgt_fysiekvoorkomen, plus_fysiekvoorkomen"
$ws.Range("F54:F68").Value = "identificatie_lokaalid, bgt_fysiekvoorkomen, plus_fysiekvoorkomen"
$ws.Range("F72").Value = "identificatie_lokaalid, bgt_fysiekvoorkomen, plus_fysiekvoorkomen"
$ws.Range("F75").Value = "identificatie_lokaalid, bgt_fysiekvoorkomen, plus_fysiekvoorkomen"
$ws.Range("F87:F94").Value = "identificatie_lokaalid, bgt_fysiekvoorkomen, plus_fysiekvoorkomen"
$ws.Range("F127").Value = "identificatie_lokaalid, bgt_fysiekvoorkomen, plus_fysiekvoorkomen"

# "lokaalid, bgtfunctie, plusfunct" -> "identificatie_lokaalid, bgt_functie, plus_functie"
$ws.Range("F96:F98").Value = "identificatie_lokaalid, bgt_functie, plus_functie"
$ws.Range("F109:F111").Value = "identificatie_lokaalid, bgt_functie, plus_functie"
$ws.Range("F116:F126").Value = "identificatie_lokaalid, bgt_functie, plus_functie"
$ws.Range("F128:F132").Value = "identificatie_lokaalid, bgt_functie, plus_functie"

# "BGT_WGL_ov-baan" -> "BGT_WGL_ov_baan"
$ws.Range("B119").Value = "BGT_WGL_ov_baan"

# ---------------------------------------------------------------------------
# 2) AutoFilter bookkeeping: add the extra sheet-scoped _FilterDatabase names
#    that show up alongside the existing _xlnm._FilterDatabase pair.
# ---------------------------------------------------------------------------
$ws.Names.Add("_xlnm._FilterDatabase_0", "=Blad1!`$A`$1:`$J`$132")
$ws.Names.Add("_xlnm._FilterDatabase_0_0", "=Blad1!`$A`$1:`$J`$132")

# ---------------------------------------------------------------------------
# 3) View state: scroll the frozen pane down and select the cells the author
#    ended up touching/looking at.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 107
$sel = $ws.Range("I10:I11,I111,I127")
$sel.Select()

# ---------------------------------------------------------------------------
# 4) Header row 2 gets a touch shorter.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 13.8

# ---------------------------------------------------------------------------
# 5) Column widths A:J grow (mostly a consistent re-wrap of the new, longer
#    text; column F grows much more since its text got a lot longer).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.020917678812367
$ws.Columns.Item(2).ColumnWidth = 43.385290148448064
$ws.Columns.Item(3).ColumnWidth = 22.818488529014868
$ws.Columns.Item(4).ColumnWidth = 21.676788124156566
$ws.Columns.Item(5).ColumnWidth = 44.640350877192965
$ws.Columns.Item(6).ColumnWidth = 61.551282051282065
$ws.Columns.Item(7).ColumnWidth = 44.069500674763866
$ws.Columns.Item(8).ColumnWidth = 24.761808367071566
$ws.Columns.Item(9).ColumnWidth = 19.275978407557368
$ws.Columns.Item(10).ColumnWidth = 19.733468286099868
